$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 is the "apple-core" entry.
# Update the density value (B2) from 0.909 to 0.96
$ws.Range("B2").Value = 0.96

# Update the Note (D2) to prepend the new source reference found during the
# fruit waste analysis, keeping the rest of the original note text intact.
$ws.Range("D2").Value = "https://files.eric.ed.gov/fulltext/EJ1107681.pdf found 0.96g/cm3;168202;`t Apples, raw, golden delicious, with skin; 1 cup is 236.588ml; 1 cup of water is 236.588g (4 degrees C); source from http://foodinfo.us/Densities.aspx; Assume core thickness is 2 cm"

# Update the selected cell to match the author's final cursor position.
$ws.Range("B11").Select()
